$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.3194395642580616
$ws.Cells.Item(2, 3).Value = 0.04439482012189444
$ws.Cells.Item(2, 4).Value = 0.1841706617995271
$ws.Cells.Item(2, 5).Value = 0.1605469713034111
$ws.Cells.Item(2, 6).Value = 1.458553068320995
$ws.Cells.Item(2, 9).Value = 0.7272070302619085
$ws.Cells.Item(2, 10).Value = 0.1822206034043532
$ws.Cells.Item(2, 11).Value = 0.3387114728725464
$ws.Cells.Item(2, 13).Value = 0.2267650527506433
$ws.Cells.Item(2, 14).Value = 1.782897933633262
$ws.Cells.Item(2, 15).Value = 3.535481513871957
$ws.Cells.Item(3, 2).Value = 0.2884028525818962
$ws.Cells.Item(3, 3).Value = 0.03943736159244793
$ws.Cells.Item(3, 4).Value = 0.1807656356292568
$ws.Cells.Item(3, 5).Value = 0.159330819151652
$ws.Cells.Item(3, 6).Value = 1.462025641079492
$ws.Cells.Item(3, 9).Value = 0.732623383212438
$ws.Cells.Item(3, 10).Value = 0.1820152255282252
$ws.Cells.Item(3, 11).Value = 0.3046738519000485
$ws.Cells.Item(3, 13).Value = 0.2163156387702116
$ws.Cells.Item(3, 14).Value = 1.800013050368399
$ws.Cells.Item(3, 15).Value = 3.552438498571661
$ws.Cells.Item(4, 2).Value = 0.2693925950910909
$ws.Cells.Item(4, 3).Value = 0.036393739878406
$ws.Cells.Item(4, 4).Value = 0.1787478435328183
$ws.Cells.Item(4, 5).Value = 0.1586532487716106
$ws.Cells.Item(4, 6).Value = 1.464836379882435
$ws.Cells.Item(4, 9).Value = 0.736274236775099
$ws.Cells.Item(4, 10).Value = 0.1819707202895415
$ws.Cells.Item(4, 11).Value = 0.2838141513817334
$ws.Cells.Item(4, 13).Value = 0.2099892080369123
$ws.Cells.Item(4, 14).Value = 1.811060716742305
$ws.Cells.Item(4, 15).Value = 3.564513812197831
$ws.Cells.Item(5, 2).Value = 0.2616579329380215
$ws.Cells.Item(5, 3).Value = 0.03515356876310705
$ws.Cells.Item(5, 4).Value = 0.1779440075164018
$ws.Cells.Item(5, 5).Value = 0.1583945698804499
$ws.Cells.Item(5, 6).Value = 1.466152597084033
$ws.Cells.Item(5, 9).Value = 0.7378438015603876
$ws.Cells.Item(5, 10).Value = 0.1819731315527093
$ws.Cells.Item(5, 11).Value = 0.2753241151065851
$ws.Cells.Item(5, 13).Value = 0.2074338411865568
$ws.Cells.Item(5, 14).Value = 1.815698279955531
$ws.Cells.Item(5, 15).Value = 3.569853222022701
$ws.Cells.Item(6, 2).Value = 0.2603743481243725
$ws.Cells.Item(6, 3).Value = 0.03494764886640667
$ws.Cells.Item(6, 4).Value = 0.1778116471869708
$ws.Cells.Item(6, 5).Value = 0.1583526709745087
$ws.Cells.Item(6, 6).Value = 1.46638147647532
$ws.Cells.Item(6, 9).Value = 0.7381093690524168
$ws.Cells.Item(6, 10).Value = 0.181974773816922
$ws.Cells.Item(6, 11).Value = 0.2739149981540976
$ws.Cells.Item(6, 13).Value = 0.2070109012825725
$ws.Cells.Item(6, 14).Value = 1.816476531706601
$ws.Cells.Item(6, 15).Value = 3.570765118965383
$ws.Cells.Item(7, 2).Value = 0.2692882327829125
$ws.Cells.Item(7, 3).Value = 0.0363770138939401
$ws.Cells.Item(7, 4).Value = 0.1787369279651614
$ws.Cells.Item(7, 5).Value = 0.1586496894742204
$ws.Cells.Item(7, 6).Value = 1.464853438962265
$ws.Cells.Item(7, 9).Value = 0.7362950731662252
$ws.Cells.Item(7, 10).Value = 0.1819706695737722
$ws.Cells.Item(7, 11).Value = 0.2836996087389991
$ws.Cells.Item(7, 13).Value = 0.2099546532866974
$ws.Cells.Item(7, 14).Value = 1.811122711660497
$ws.Cells.Item(7, 15).Value = 3.564584125991175
$ws.Cells.Item(8, 2).Value = 0.3087287655411899
$ws.Cells.Item(8, 3).Value = 0.04268546499628201
$ws.Cells.Item(8, 4).Value = 0.1829815338158767
$ws.Cells.Item(8, 5).Value = 0.1601133174701808
$ws.Cells.Item(8, 6).Value = 1.459609685274124
$ws.Cells.Item(8, 9).Value = 0.7290071194891006
$ws.Cells.Item(8, 10).Value = 0.1821328705249456
$ws.Cells.Item(8, 11).Value = 0.3269673964486515
$ws.Cells.Item(8, 13).Value = 0.2231436166929797
$ws.Cells.Item(8, 14).Value = 1.788687396432916
$ws.Cells.Item(8, 15).Value = 3.540983155448288
$ws.Cells.Item(9, 2).Value = 0.3864215199585033
$ws.Cells.Item(9, 3).Value = 0.05505650435607379
$ws.Cells.Item(9, 4).Value = 0.1918799720279907
$ws.Cells.Item(9, 5).Value = 0.1635303835630531
$ws.Cells.Item(9, 6).Value = 1.45470374531218
$ws.Cells.Item(9, 9).Value = 0.7172940576753675
$ws.Cells.Item(9, 10).Value = 0.1830973222097398
$ws.Cells.Item(9, 11).Value = 0.4121104515646721
$ws.Cells.Item(9, 13).Value = 0.2497109429754687
$ws.Cells.Item(9, 14).Value = 1.74896564791001
$ws.Cells.Item(9, 15).Value = 3.507891647955432
$ws.Cells.Item(10, 2).Value = 0.4436968565299821
$ws.Cells.Item(10, 3).Value = 0.06414379003196302
$ws.Cells.Item(10, 4).Value = 0.1987639042002201
$ws.Cells.Item(10, 5).Value = 0.166372340129552
$ws.Cells.Item(10, 6).Value = 1.454369587396641
$ws.Cells.Item(10, 9).Value = 0.7102584158738985
$ws.Cells.Item(10, 10).Value = 0.1841988407840986
$ws.Cells.Item(10, 11).Value = 0.4748260051609918
$ws.Cells.Item(10, 13).Value = 0.2696520957710078
$ws.Cells.Item(10, 14).Value = 1.722383528284893
$ws.Cells.Item(10, 15).Value = 3.49160906164559
$ws.Cells.Item(11, 2).Value = 0.4697914353129988
$ws.Cells.Item(11, 3).Value = 0.06827712184980328
$ws.Cells.Item(11, 4).Value = 0.2019699361011646
$ws.Cells.Item(11, 5).Value = 0.1677367983865565
$ws.Cells.Item(11, 6).Value = 1.454926148821798
$ws.Cells.Item(11, 9).Value = 0.7073982118876181
$ws.Cells.Item(11, 10).Value = 0.1847850570588179
$ws.Cells.Item(11, 11).Value = 0.5033882239912373
$ws.Cells.Item(11, 13).Value = 0.2788141604497412
$ws.Cells.Item(11, 14).Value = 1.710854663553303
$ws.Cells.Item(11, 15).Value = 3.485943132268517
$ws.Cells.Item(12, 2).Value = 0.4796780476887932
$ws.Cells.Item(12, 3).Value = 0.06984218311505686
$ws.Cells.Item(12, 4).Value = 0.2031946012229184
$ws.Cells.Item(12, 5).Value = 0.1682637425173645
$ws.Cells.Item(12, 6).Value = 1.455238639386209
$ws.Cells.Item(12, 9).Value = 0.7063640362688304
$ws.Cells.Item(12, 10).Value = 0.1850192592821145
$ws.Cells.Item(12, 11).Value = 0.514208229730599
$ws.Cells.Item(12, 13).Value = 0.2822964811779514
$ws.Cells.Item(12, 14).Value = 1.706569980049891
$ws.Cells.Item(12, 15).Value = 3.484047744935879
$ws.Cells.Item(13, 2).Value = 0.4775485680616498
$ws.Cells.Item(13, 3).Value = 0.06950512654360352
$ws.Cells.Item(13, 4).Value = 0.2029303774024527
$ws.Cells.Item(13, 5).Value = 0.1681498004144331
$ws.Cells.Item(13, 6).Value = 1.455166816891534
$ws.Cells.Item(13, 9).Value = 0.7065845890042262
$ws.Cells.Item(13, 10).Value = 0.1849682768113823
$ws.Cells.Item(13, 11).Value = 0.5118777741094505
$ws.Cells.Item(13, 13).Value = 0.2815459330053045
$ws.Cells.Item(13, 14).Value = 1.707489158351761
$ws.Cells.Item(13, 15).Value = 3.484444826859772
$ws.Cells.Item(14, 2).Value = 0.470604712584219
$ws.Cells.Item(14, 3).Value = 0.06840588360722677
$ws.Cells.Item(14, 4).Value = 0.2020704779744591
$ws.Cells.Item(14, 5).Value = 0.1677799451796069
$ws.Cells.Item(14, 6).Value = 1.454949819642763
$ws.Cells.Item(14, 9).Value = 0.7073121492882457
$ws.Cells.Item(14, 10).Value = 0.1848040803664475
$ws.Cells.Item(14, 11).Value = 0.5042783132300599
$ws.Cells.Item(14, 13).Value = 0.2791003967827876
$ws.Cells.Item(14, 14).Value = 1.71050053562765
$ws.Cells.Item(14, 15).Value = 3.485782184784739
$ws.Cells.Item(15, 2).Value = 0.4663520563918269
$ws.Cells.Item(15, 3).Value = 0.06773254511311677
$ws.Cells.Item(15, 4).Value = 0.2015451435021305
$ws.Cells.Item(15, 5).Value = 0.1675547320362369
$ws.Cells.Item(15, 6).Value = 1.454830146457326
$ws.Cells.Item(15, 9).Value = 0.7077641714203509
$ws.Cells.Item(15, 10).Value = 0.1847050952407514
$ws.Cells.Item(15, 11).Value = 0.4996239426338605
$ws.Cells.Item(15, 13).Value = 0.2776041021230427
$ws.Cells.Item(15, 14).Value = 1.712355646565843
$ws.Cells.Item(15, 15).Value = 3.486633929905139
$ws.Cells.Item(16, 2).Value = 0.4419922668679988
$ws.Cells.Item(16, 3).Value = 0.06387365088511388
$ws.Cells.Item(16, 4).Value = 0.1985558728129604
$ws.Cells.Item(16, 5).Value = 0.1662846070079844
$ws.Cells.Item(16, 6).Value = 1.454347459495793
$ws.Cells.Item(16, 9).Value = 0.7104521838597577
$ws.Cells.Item(16, 10).Value = 0.184162240920493
$ws.Cells.Item(16, 11).Value = 0.4729600038813544
$ws.Cells.Item(16, 13).Value = 0.2690551418616565
$ws.Cells.Item(16, 14).Value = 1.723148306947505
$ws.Cells.Item(16, 15).Value = 3.492014362532643
$ws.Cells.Item(17, 2).Value = 0.427058110318228
$ws.Cells.Item(17, 3).Value = 0.0615061634314884
$ws.Cells.Item(17, 4).Value = 0.1967410598429638
$ws.Cells.Item(17, 5).Value = 0.1655237404939811
$ws.Cells.Item(17, 6).Value = 1.454232715915367
$ws.Cells.Item(17, 9).Value = 0.7121883489067251
$ws.Cells.Item(17, 10).Value = 0.1838510022023669
$ws.Cells.Item(17, 11).Value = 0.4566104820271164
$ws.Cells.Item(17, 13).Value = 0.2638337274120985
$ws.Cells.Item(17, 14).Value = 1.72991358765178
$ws.Cells.Item(17, 15).Value = 3.495760885537067
$ws.Cells.Item(18, 2).Value = 0.418472152152674
$ws.Cells.Item(18, 3).Value = 0.06014440233161622
$ws.Cells.Item(18, 4).Value = 0.1957042455129852
$ws.Cells.Item(18, 5).Value = 0.1650928559844118
$ws.Cells.Item(18, 6).Value = 1.454233422070523
$ws.Cells.Item(18, 9).Value = 0.713218981146337
$ws.Cells.Item(18, 10).Value = 0.1836800003621235
$ws.Cells.Item(18, 11).Value = 0.4472097722238004
$ws.Cells.Item(18, 13).Value = 0.2608390603497455
$ws.Cells.Item(18, 14).Value = 1.733857833716041
$ws.Cells.Item(18, 15).Value = 3.498079678553097
$ws.Cells.Item(19, 2).Value = 0.4155657592212378
$ws.Cells.Item(19, 3).Value = 0.05968332774365592
$ws.Cells.Item(19, 4).Value = 0.1953544064676009
$ws.Cells.Item(19, 5).Value = 0.1649481259500511
$ws.Cells.Item(19, 6).Value = 1.454245122411535
$ws.Cells.Item(19, 9).Value = 0.7135734383317249
$ws.Cells.Item(19, 10).Value = 0.1836234793779781
$ws.Cells.Item(19, 11).Value = 0.4440274037990832
$ws.Cells.Item(19, 13).Value = 0.2598265925674781
$ws.Cells.Item(19, 14).Value = 1.73520239598832
$ws.Cells.Item(19, 15).Value = 3.498892936468678
$ws.Cells.Item(20, 2).Value = 0.4286474902186228
$ws.Cells.Item(20, 3).Value = 0.06175819164961638
$ws.Cells.Item(20, 4).Value = 0.1969335241752788
$ws.Cells.Item(20, 5).Value = 0.1656040382178006
$ws.Cells.Item(20, 6).Value = 1.454238028139628
$ws.Cells.Item(20, 9).Value = 0.7120002157072172
$ws.Cells.Item(20, 10).Value = 0.1838833048243984
$ws.Cells.Item(20, 11).Value = 0.4583505995758514
$ws.Cells.Item(20, 13).Value = 0.2643886722104511
$ws.Cells.Item(20, 14).Value = 1.729187922851199
$ws.Cells.Item(20, 15).Value = 3.495345100933434
$ws.Cells.Item(21, 2).Value = 0.472644155581122
$ws.Cells.Item(21, 3).Value = 0.06872876228698033
$ws.Cells.Item(21, 4).Value = 0.2023227641316225
$ws.Cells.Item(21, 5).Value = 0.1678883027366389
$ws.Cells.Item(21, 6).Value = 1.455010797263313
$ws.Cells.Item(21, 9).Value = 0.7070971196353888
$ws.Cells.Item(21, 10).Value = 0.1848519775486963
$ws.Cells.Item(21, 11).Value = 0.5065103520493608
$ws.Cells.Item(21, 13).Value = 0.2798183630344013
$ws.Cells.Item(21, 14).Value = 1.709613821348324
$ws.Cells.Item(21, 15).Value = 3.485382582131763
$ws.Cells.Item(22, 2).Value = 0.5014283491991307
$ws.Cells.Item(22, 3).Value = 0.07328357554058584
$ws.Cells.Item(22, 4).Value = 0.2059067421415364
$ws.Cells.Item(22, 5).Value = 0.1694409439306241
$ws.Cells.Item(22, 6).Value = 1.45610874933692
$ws.Cells.Item(22, 9).Value = 0.7041778061416899
$ws.Cells.Item(22, 10).Value = 0.1855562461516342
$ws.Cells.Item(22, 11).Value = 0.538009243054006
$ws.Cells.Item(22, 13).Value = 0.2899773384331255
$ws.Cells.Item(22, 14).Value = 1.697293431696187
$ws.Cells.Item(22, 15).Value = 3.48032963523724
$ws.Cells.Item(23, 2).Value = 0.4860631367975827
$ws.Cells.Item(23, 3).Value = 0.07085268662369515
$ws.Cells.Item(23, 4).Value = 0.2039882848689842
$ws.Cells.Item(23, 5).Value = 0.1686068195907389
$ws.Cells.Item(23, 6).Value = 1.455468553610487
$ws.Cells.Item(23, 9).Value = 0.7057098145275482
$ws.Cells.Item(23, 10).Value = 0.185173860171524
$ws.Cells.Item(23, 11).Value = 0.5211957127045537
$ws.Cells.Item(23, 13).Value = 0.2845485264382717
$ws.Cells.Item(23, 14).Value = 1.703825821957159
$ws.Cells.Item(23, 15).Value = 3.482893131832355
$ws.Cells.Item(24, 2).Value = 0.4279289318751864
$ws.Cells.Item(24, 3).Value = 0.06164425174588928
$ws.Cells.Item(24, 4).Value = 0.1968464906540532
$ws.Cells.Item(24, 5).Value = 0.1655677152179535
$ws.Cells.Item(24, 6).Value = 1.454235418787391
$ws.Cells.Item(24, 9).Value = 0.7120851695276897
$ws.Cells.Item(24, 10).Value = 0.183868676095166
$ws.Cells.Item(24, 11).Value = 0.4575638959768469
$ws.Cells.Item(24, 13).Value = 0.2641377592448748
$ws.Cells.Item(24, 14).Value = 1.729515825221105
$ws.Cells.Item(24, 15).Value = 3.495532563587517
$ws.Cells.Item(25, 2).Value = 0.3653680542051063
$ws.Cells.Item(25, 3).Value = 0.05170997943602629
$ws.Cells.Item(25, 4).Value = 0.1894116032949
$ws.Cells.Item(25, 5).Value = 0.1625476303474507
$ws.Cells.Item(25, 6).Value = 1.455456083501055
$ws.Cells.Item(25, 9).Value = 0.7201869061991992
$ws.Cells.Item(25, 10).Value = 0.1827673282248909
$ws.Cells.Item(25, 11).Value = 0.3890473657480129
$ws.Cells.Item(25, 13).Value = 0.2424491439794423
$ws.Cells.Item(25, 14).Value = 1.75925431276824
$ws.Cells.Item(25, 15).Value = 3.515432863162175
